$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated "K" (column G) values from std/mean recalculation (s_vals)
$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 2
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 2
$ws.Range("G9").Value = 2
$ws.Range("G10").Value = 1
